# Generate Report for Handback
# Regenerates the handback-status report: the first test file's GUID changed
# from 6a81eca0-262a-4276-91c6-6afdc564ae7b -> 57f738f8-fada-4700-8649-83148c0867a0
# and the second test file's GUID changed from
# b8051b9e-404c-4aa7-acd6-ac3064f8eac5 -> ffff2090f858-e132-4bdd-ae11-b9874b967283,
# along with refreshed handoff/handback file names and timestamps.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # Overview
$ws2 = $wb.Worksheets.Item(2)   # zh-cn
$ws3 = $wb.Worksheets.Item(3)   # de-de

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ws1.Range("A2").Value = "57f738f8-fada-4700-8649-83148c0867a0.md"
$ws1.Range("B2").Value = "e2e\57f738f8-fada-4700-8649-83148c0867a0.md"
$ws1.Range("G2").Value = "2016-08-27 23:01:36"

$ws1.Range("A3").Value = "ffff2090f858-e132-4bdd-ae11-b9874b967283.md"
$ws1.Range("B3").Value = "e2e\ffff2090f858-e132-4bdd-ae11-b9874b967283.md"
$ws1.Range("G3").Value = "2016-08-27 23:01:36"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws2.Range("A2").Value = "57f738f8-fada-4700-8649-83148c0867a0.md"
$ws2.Range("G2").Value = "57f738f8-fada-4700-8649-83148c0867a0.79a5a110f43a5efdcd6b1813773f0934f9958326.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-27 23:01:32"
$ws2.Range("I2").Value = "57f738f8-fada-4700-8649-83148c0867a0.md"
$ws2.Range("J2").Value = "57f738f8-fada-4700-8649-83148c0867a0.79a5a110f43a5efdcd6b1813773f0934f9958326.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-08-27 23:01:48"

$ws2.Range("A3").Value = "ffff2090f858-e132-4bdd-ae11-b9874b967283.md"
$ws2.Range("G3").Value = "57f738f8-fada-4700-8649-83148c0867a0.79a5a110f43a5efdcd6b1813773f0934f9958326.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-27 23:01:32"
$ws2.Range("I3").Value = "ffff2090f858-e132-4bdd-ae11-b9874b967283.md"
$ws2.Range("J3").Value = "57f738f8-fada-4700-8649-83148c0867a0.79a5a110f43a5efdcd6b1813773f0934f9958326.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-08-27 23:01:48"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws3.Range("A2").Value = "57f738f8-fada-4700-8649-83148c0867a0.md"
$ws3.Range("G2").Value = "57f738f8-fada-4700-8649-83148c0867a0.79a5a110f43a5efdcd6b1813773f0934f9958326.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-27 23:01:36"
$ws3.Range("I2").Value = "57f738f8-fada-4700-8649-83148c0867a0.md"
$ws3.Range("J2").Value = "57f738f8-fada-4700-8649-83148c0867a0.79a5a110f43a5efdcd6b1813773f0934f9958326.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-27 23:01:55"

$ws3.Range("A3").Value = "ffff2090f858-e132-4bdd-ae11-b9874b967283.md"
$ws3.Range("G3").Value = "57f738f8-fada-4700-8649-83148c0867a0.79a5a110f43a5efdcd6b1813773f0934f9958326.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-27 23:01:36"
$ws3.Range("I3").Value = "ffff2090f858-e132-4bdd-ae11-b9874b967283.md"
$ws3.Range("J3").Value = "57f738f8-fada-4700-8649-83148c0867a0.79a5a110f43a5efdcd6b1813773f0934f9958326.de-de.xlf"
$ws3.Range("K3").Value = "2016-08-27 23:01:55"
